$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing data shifts down to rows 2-11
$ws.Rows.Item(1).Insert()

# Populate the header row
$ws.Range("A1").Value = "CNE"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "DateofBirth"
$ws.Range("E1").Value = "ClasseName"
$ws.Range("F1").Value = "Phone"
$ws.Range("G1").Value = "Email"

# Update the first data row's CNE value (was 19000051 -> 18000041); the
# dependent formula chain (A3..A11 each "=prev+1") recalculates automatically
$ws.Range("A2").Value = 18000041

$ws.Range("F7").Select()
